$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7
$ws.Cells.Item($row, 1).Value = 42611.887013888889
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item($row, 2).Value = 16
$ws.Cells.Item($row, 3).Value = 60
$ws.Cells.Item($row, 4).Value = 37
$ws.Cells.Item($row, 5).Value = 70
$ws.Cells.Item($row, 6).Value = 30
$ws.Cells.Item($row, 7).Value = 17002
$ws.Cells.Item($row, 8).Value = 16209
$ws.Cells.Item($row, 9).Value = 2662
$ws.Cells.Item($row, 10).Value = 336
$ws.Cells.Item($row, 11).Value = 211
$ws.Cells.Item($row, 12).Value = 14
$ws.Cells.Item($row, 13).Value = 6
$ws.Cells.Item($row, 14).Value = "Noun"
